$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = "Josh Hart"
$ws.Range("B15").Value = "SG,SF,PF"
$ws.Range("C15").Value = "New York Knicks"

$ws.Range("A16").Value = "Guerschon Yabusele"
$ws.Range("B16").Value = "PF,C"
$ws.Range("C16").Value = "Philadelphia 76ers"

$ws.Range("A19").Value = "De'Andre Hunter"
$ws.Range("B19").Value = "SF,PF"
$ws.Range("C19").Value = "Atlanta Hawks"
